# Update "想去人数" (want-to-go count) figures on the 展览 and 全部类型 sheets
# to reflect the latest scrape (gh-pages output generated at 456a3b4).

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 1068
$ws1.Range("F4").Value = 1582
$ws1.Range("F5").Value = 727
$ws1.Range("F6").Value = 33

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 1068
$ws4.Range("F4").Value = 1582
$ws4.Range("F6").Value = 727
$ws4.Range("F7").Value = 33
